$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 250555.8564151394
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 251711.4028993678

$ws.Range("B3").Value = 0.1190320826869504
$ws.Range("C3").Value = 117.745847958593
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 1254.439558216954
